$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme

Write-Host "Theme props probe:"
$theme | Get-Member | ForEach-Object { Write-Host $_.Name }
